$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the window position (xWindow 280 -> 120)
$wb.Windows.Item(1).Left = 120
$excel.Left = 120

# 2. Insert three new rows before current row 112 (shift rows 112-113 down to 115-116)
$ws.Rows.Item(112).Insert()
$ws.Rows.Item(112).Insert()
$ws.Rows.Item(112).Insert()

# 3. Fill the three new rows (112, 113, 114)
$refText = "Subrahmanyam & Reddy 2008. — von Hinüber 2013."

# Row 112
$ws.Cells.Item(112, 1).Value = 109
$ws.Cells.Item(112, 2).Value = "—"
$ws.Cells.Item(112, 4).Value = "—"
$ws.Cells.Item(112, 6).Value = 6
$ws.Cells.Item(112, 8).Value = "MIA"
$ws.Cells.Item(112, 9).Value = "Phanigiri"
$ws.Cells.Item(112, 13).Value = $refText

# Row 113
$ws.Cells.Item(113, 1).Value = 110
$ws.Cells.Item(113, 2).Value = "—"
$ws.Cells.Item(113, 4).Value = "—"
$ws.Cells.Item(113, 6).Value = 4
$ws.Cells.Item(113, 8).Value = "MIA"
$ws.Cells.Item(113, 9).Value = "Phanigiri"
$ws.Cells.Item(113, 13).Value = $refText

# Row 114
$ws.Cells.Item(114, 1).Value = 111
$ws.Cells.Item(114, 2).Value = "—"
$ws.Cells.Item(114, 4).Value = "—"
$ws.Cells.Item(114, 6).Value = 3
$ws.Cells.Item(114, 8).Value = "MIA"
$ws.Cells.Item(114, 9).Value = "Phanigiri"
$ws.Cells.Item(114, 13).Value = $refText

# Set row heights explicitly
$ws.Rows.Item(112).RowHeight = 32
$ws.Rows.Item(113).RowHeight = 32
$ws.Rows.Item(114).RowHeight = 32

# Apply the "red reference" style (style index 4 = red Calibri font) to the
# References column in the three new rows, matching M111's new style.
$ws.Cells.Item(112, 13).Font.Color = 255  # RGB(255,0,0)
$ws.Cells.Item(113, 13).Font.Color = 255
$ws.Cells.Item(114, 13).Font.Color = 255

# 4. Update M111 style to red font (style index 4) and keep value/text the same
$ws.Cells.Item(111, 13).Font.Color = 255  # RGB(255,0,0) = FF0000 -> represented as BGR in VBA color = 255

# Clear the stray formatting that Insert() propagated into column E of the
# new rows (the source rows had no E value there).
$ws.Cells.Item(112, 5).Clear()
$ws.Cells.Item(113, 5).Clear()
$ws.Cells.Item(114, 5).Clear()

# 5. Update the selection to M111
$ws.Range("M111").Select()
